# Atualização automática de MARAU.xlsx
#
# - "Paineis DARQ" -> "PAINEIS DARQ"
# - "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Remove the "Desarquivamentos Pendentes" sheet (no longer needed; its two
#   unique strings and dedicated cell styles go away with it)
# - "DGC" sheet stays as-is (content untouched)

$wb = $excel.ActiveWorkbook

$painelSheet = $wb.Worksheets.Item("Paineis DARQ")
$painelSheet.Name = "PAINEIS DARQ"

$recolhSheet = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$recolhSheet.Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$pendentesSheet = $wb.Worksheets.Item("Desarquivamentos Pendentes")
[void]$pendentesSheet.Delete()
$excel.DisplayAlerts = $true
